# Swap the order of "Recorded By" names in column G for rows where the
# value is exactly "<name>, System" -> "System, <name>".
#
# This mirrors the upstream diff where every G-column cell whose text was
# "<something>, System" became "System, <something>" (the "System" token
# was moved to the front of the comma-separated list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    $text = [string]$val
    if ($text -eq "") { continue }

    $commaIndex = $text.IndexOf(",")
    if ($commaIndex -lt 0) { continue }

    $before = $text.Substring(0, $commaIndex)
    $after = $text.Substring($commaIndex + 1).Trim()

    if ($after -eq "System") {
        $cell.Value2 = "System, " + $before
    }
}
